$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [string][char]0x2083

$ws.Range("D2").Value = "'63.056.73"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.15%  "
$ws.Range("E2").NumberFormat = "General"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.029.80"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.91%  "
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").NumberFormat = "General"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'596.19"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.58%  "
$ws.Range("E5").NumberFormat = "General"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'153.09"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +7.34%  "
$ws.Range("E6").NumberFormat = "General"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.09%  "
$ws.Range("E7").NumberFormat = "General"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.028.20"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.95%  "
$ws.Range("E8").NumberFormat = "General"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.514"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.31%  "
$ws.Range("E9").NumberFormat = "General"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'7.02"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +16.66%  "
$ws.Range("E10").NumberFormat = "General"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +2.22%  "
$ws.Range("E11").NumberFormat = "General"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.464"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +3.08%  "
$ws.Range("E12").NumberFormat = "General"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000233"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.41%  "
$ws.Range("E13").NumberFormat = "General"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'35.77"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +5.31%  "
$ws.Range("E14").NumberFormat = "General"
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.10%  "
$ws.Range("E15").NumberFormat = "General"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.528.11"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.85%  "
$ws.Range("E16").NumberFormat = "General"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'Polkadot"
$ws.Range("B17").NumberFormat = "General"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C17").NumberFormat = "General"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'7.08"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.99%  "
$ws.Range("E17").NumberFormat = "General"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'WrappedBTC"
$ws.Range("B18").NumberFormat = "General"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C18").NumberFormat = "General"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'62.968.46"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'3.025.12"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.86%  "
$ws.Range("E19").NumberFormat = "General"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'449.87"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.81%  "
$ws.Range("E20").NumberFormat = "General"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.25"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.29%  "
$ws.Range("E21").NumberFormat = "General"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.697"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.61%  "
$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.54"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.29%  "
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'RenderToken"
$ws.Range("B24").NumberFormat = "General"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C24").NumberFormat = "General"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'11.45"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +9.30%  "
$ws.Range("E24").NumberFormat = "General"
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'Litecoin"
$ws.Range("B25").NumberFormat = "General"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C25").NumberFormat = "General"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'82.92"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.85%  "
$ws.Range("E25").NumberFormat = "General"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.31"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +6.51%  "
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.37"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.57%  "
$ws.Range("E27").NumberFormat = "General"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.00%  "
$ws.Range("E28").NumberFormat = "General"
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.51"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.39%  "
$ws.Range("E29").NumberFormat = "General"
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.28"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +11.52%  "
$ws.Range("E30").NumberFormat = "General"
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.70"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.22%  "
$ws.Range("E31").NumberFormat = "General"
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.03%  "
$ws.Range("E32").NumberFormat = "General"
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'27.74"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.63%  "
$ws.Range("E33").NumberFormat = "General"
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.111"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.96%  "
$ws.Range("E34").NumberFormat = "General"
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0" + $sub3 + "0877"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +8.70%  "
$ws.Range("E35").NumberFormat = "General"
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +3.29%  "
$ws.Range("E36").NumberFormat = "General"
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.89"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.64%  "
$ws.Range("E37").NumberFormat = "General"
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'3.14"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +11.40%  "
$ws.Range("E38").NumberFormat = "General"
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +8.45%  "
$ws.Range("E39").NumberFormat = "General"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.10"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +4.27%  "
$ws.Range("E40").NumberFormat = "General"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'50.57"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.15%  "
$ws.Range("E41").NumberFormat = "General"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'9.04"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'44.67"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +17.13%  "
$ws.Range("E43").NumberFormat = "General"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.307"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +15.28%  "
$ws.Range("E44").NumberFormat = "General"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'391.24"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.11%  "
$ws.Range("E45").NumberFormat = "General"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0361"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +3.67%  "
$ws.Range("E46").NumberFormat = "General"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.709.49"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.17%  "
$ws.Range("E47").NumberFormat = "General"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'133.91"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.77%  "
$ws.Range("E48").NumberFormat = "General"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'26.83"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +15.65%  "
$ws.Range("E49").NumberFormat = "General"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.02%  "
$ws.Range("E50").NumberFormat = "General"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.28"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +6.65%  "
$ws.Range("E51").NumberFormat = "General"
$ws.Range("E51").Style = "Normal"
